$wb = $excel.ActiveWorkbook

# --- Sheet "Daily Orders": insert a new order row at the top of the data (row 2) ---
$ws1 = $wb.Worksheets.Item("Daily Orders")
$ws1.Rows("2:2").Insert()

$ws1.Range("A2").Value = 11
$ws1.Range("B2").Value = "2026-01-13 22:38"
$ws1.Range("C2").Value = "Phantom"

# D2 ("420") and J2 ("2026-01-14") look numeric/date-like, so force text
# formatting before the write and then restore the Normal style so no
# stray formatting is left behind on the cell.
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "420"
$ws1.Range("D2").Style = "Normal"

$ws1.Range("F2").Value = "Upma x1"
$ws1.Range("G2").Value = 30
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"

$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "2026-01-14"
$ws1.Range("J2").Style = "Normal"

$ws1.Range("K2").Value = "15:38"
# E2, L2, M2, N2 stay blank for the new order (matching the template's
# empty-cell pattern for Phone / Notes / Cancel Reason / Feedback).

# --- Sheet "Summary": bump the aggregate counters for the new order ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 11
$ws2.Range("B2").Value = 9
$ws2.Range("G2").Value = 275

# --- Sheet "Items Breakdown": insert the new "Upma" item line ---
$ws3 = $wb.Worksheets.Item("Items Breakdown")
$ws3.Rows("4:4").Insert()
$ws3.Range("A4").Value = "Upma"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 30
